# Workbook / worksheet handles
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows right before the current row 344 (which shifts the
# existing rows 344-364 down to 347-367) so we can populate the 3 new
# weekly records while keeping all the later rows (now 347-367) intact
# with their original values.
$ws.Rows.Item(344).Insert()
$ws.Rows.Item(345).Insert()
$ws.Rows.Item(346).Insert()

# Columns A,B,C,E,F,G,H,Q,R are constant for this whole data block
# (Mercado ID, Mercado, Región, Mercado ID nro, Producto ID, Producto,
# Variedad, Unidad flag, Categoria). Re-use them for the 3 new rows.
$commonA = 9
$commonB = "Vega Central Mapocho de Santiago"
$commonC = "Metropolitana"
$commonE = 13
$commonF = 100112028
$commonG = "Sandia"
$commonH = "Sin especificar"
$commonQ = 1
$commonR = "Hortaliza"

# New row 344
$ws.Cells.Item(344, 1).Value = $commonA
$ws.Cells.Item(344, 2).Value = $commonB
$ws.Cells.Item(344, 3).Value = $commonC
$ws.Cells.Item(344, 4).Value = 44585
$ws.Cells.Item(344, 5).Value = $commonE
$ws.Cells.Item(344, 6).Value = $commonF
$ws.Cells.Item(344, 7).Value = $commonG
$ws.Cells.Item(344, 8).Value = $commonH
$ws.Cells.Item(344, 9).Value = "Extra"
$ws.Cells.Item(344, 10).Value = 97
$ws.Cells.Item(344, 11).Value = 2800
$ws.Cells.Item(344, 12).Value = 3000
$ws.Cells.Item(344, 13).Value = 2899
$ws.Cells.Item(344, 14).Value = "$/unidad"
$ws.Cells.Item(344, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(344, 16).Value = 2899
$ws.Cells.Item(344, 17).Value = $commonQ
$ws.Cells.Item(344, 18).Value = $commonR

# New row 345
$ws.Cells.Item(345, 1).Value = $commonA
$ws.Cells.Item(345, 2).Value = $commonB
$ws.Cells.Item(345, 3).Value = $commonC
$ws.Cells.Item(345, 4).Value = 44585
$ws.Cells.Item(345, 5).Value = $commonE
$ws.Cells.Item(345, 6).Value = $commonF
$ws.Cells.Item(345, 7).Value = $commonG
$ws.Cells.Item(345, 8).Value = $commonH
$ws.Cells.Item(345, 9).Value = "Primera"
$ws.Cells.Item(345, 10).Value = 160
$ws.Cells.Item(345, 11).Value = 2300
$ws.Cells.Item(345, 12).Value = 2500
$ws.Cells.Item(345, 13).Value = 2400
$ws.Cells.Item(345, 14).Value = "$/unidad"
$ws.Cells.Item(345, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(345, 16).Value = 2400
$ws.Cells.Item(345, 17).Value = $commonQ
$ws.Cells.Item(345, 18).Value = $commonR

# New row 346
$ws.Cells.Item(346, 1).Value = $commonA
$ws.Cells.Item(346, 2).Value = $commonB
$ws.Cells.Item(346, 3).Value = $commonC
$ws.Cells.Item(346, 4).Value = 44585
$ws.Cells.Item(346, 5).Value = $commonE
$ws.Cells.Item(346, 6).Value = $commonF
$ws.Cells.Item(346, 7).Value = $commonG
$ws.Cells.Item(346, 8).Value = $commonH
$ws.Cells.Item(346, 9).Value = "Segunda"
$ws.Cells.Item(346, 10).Value = 61
$ws.Cells.Item(346, 11).Value = 1800
$ws.Cells.Item(346, 12).Value = 2000
$ws.Cells.Item(346, 13).Value = 1902
$ws.Cells.Item(346, 14).Value = "$/unidad"
$ws.Cells.Item(346, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(346, 16).Value = 1902
$ws.Cells.Item(346, 17).Value = $commonQ
$ws.Cells.Item(346, 18).Value = $commonR

# Make sure the new D-column cells keep the same date/time number format
# used by the rest of the column (style index 2 in the original sheet).
$ws.Range("D344:D346").NumberFormat = "YYYY-MM-DD HH:MM:SS"
